# V 0.46-B43 (PreRel) edit: add a "Free Text" field column to the Tabelle2
# "engine merge" lookup sheet.
#
# The sheet has one data column per engine-config attribute; a header row
# (row 1, style 11) with the attribute name as a shared string, and 40 data
# rows (rows 2-40, style 12) each containing "|" placeholders, ending with
# the model-name lookup formula column.
#
# We insert one new column right before the previous "END_OF_COL" marker
# column (DP), shifting END_OF_COL and the Title/formula column one place
# to the right, and populate the new column with header "TXT" and "|" filler
# values for the data rows - i.e. a new "Free Text" attribute column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# Insert a new blank column at DP (shifts old DP->DQ, old DQ->DR, etc.)
$ws.Columns("DP").Insert()

# New header cell (row 1) gets the "TXT" label
$ws.Range("DP1").Value = "TXT"

# New data cells (rows 2-40) get the "|" placeholder used throughout the sheet
$ws.Range("DP2:DP40").Value = "|"

# Restore/update the selected cell as recorded for this sheet view
$ws.Range("CR26").Select()
